$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Min/Max threshold values for the remaining parameter rows ---
# alpha_distance_range
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 9.1

# beta_distance_range
$ws.Range("B3").Value = 4.2
$ws.Range("C3").Value = 7.1

# ratio_threshold_range
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1.5

# --- Remove the "theta_threshold_range" row entirely ---
# This shifts the old row 6 ("pie_threshold_range") up to become row 5.
$ws.Rows("5").Delete()

# The shifted-up row inherited the old row 6 formatting for column B
# (a leftover "Times New Roman" style that is no longer needed anywhere
# in the sheet). Re-normalize it to the same look as the rest of the
# data cells by copying the format from a neighboring, correctly styled
# cell instead of building a brand new style.
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New values for the (now) last data row, pie_threshold_range
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# --- Leftover UI selection state recorded in the worksheet ---
$ws.Range("D6:D7").Select()
